# Weekly update: a new week of "Apio" prices (Terminal La Palmera de La
# Serena) is published. The existing history (rows 433:544) is pushed down
# by two rows to make room for the two new rows (Primera / Segunda quality)
# at the top of that block, dated 2023-01-02 (serial 44932).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing block of rows down by 2 rows (433:544 -> 435:546).
$ws.Range("A433:A434").EntireRow.Insert()

# New row 433 - "Primera" quality for the new week.
$ws.Range("A433").Value = 8
$ws.Range("B433").Value = "Terminal La Palmera de La Serena"
$ws.Range("C433").Value = "Coquimbo"
$ws.Range("D433").Value = 44932
$ws.Range("E433").Value = 4
$ws.Range("F433").Value = 100112017
$ws.Range("G433").Value = "Apio"
$ws.Range("H433").Value = "Americana (o)"
$ws.Range("I433").Value = "Primera"
$ws.Range("J433").Value = 2000
$ws.Range("K433").Value = 8000
$ws.Range("L433").Value = 9000
$ws.Range("M433").Value = 8500
$ws.Range("N433").Value = "`$/docena de matas"
$ws.Range("O433").Value = "Provincia del Elquí"
$ws.Range("P433").Value = 1417
$ws.Range("Q433").Value = 6
$ws.Range("R433").Value = "Hortaliza"

# New row 434 - "Segunda" quality for the new week.
$ws.Range("A434").Value = 8
$ws.Range("B434").Value = "Terminal La Palmera de La Serena"
$ws.Range("C434").Value = "Coquimbo"
$ws.Range("D434").Value = 44932
$ws.Range("E434").Value = 4
$ws.Range("F434").Value = 100112017
$ws.Range("G434").Value = "Apio"
$ws.Range("H434").Value = "Americana (o)"
$ws.Range("I434").Value = "Segunda"
$ws.Range("J434").Value = 1300
$ws.Range("K434").Value = 6000
$ws.Range("L434").Value = 7000
$ws.Range("M434").Value = 6500
$ws.Range("N434").Value = "`$/docena de matas"
$ws.Range("O434").Value = "Provincia del Elquí"
$ws.Range("P434").Value = 1083
$ws.Range("Q434").Value = 6
$ws.Range("R434").Value = "Hortaliza"
